# "Add file in admin branch"
#
# The original document is a single paragraph containing the run "dfdf".
# This change:
#   1. Splits that run into two runs ("D" + "fdf") and wraps them with
#      proofErr spell-check markers (spellStart/spellEnd) - as Word does
#      when it re-flags a word after an edit splits a run.
#   2. Inserts a new paragraph after it containing the text "Admin2",
#      which takes over the trailing "_GoBack" bookmark that Word
#      maintains at the last edit position.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r = $p1.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:proofErr w:type="spellStart"/>' +
              '<w:r><w:t>D</w:t></w:r>' +
              '<w:r><w:t>fdf</w:t></w:r>' +
              '<w:proofErr w:type="spellEnd"/>' +
            '</w:p>' +
            '<w:p>' +
              '<w:r><w:t>Admin2</w:t></w:r>' +
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
              '<w:bookmarkEnd w:id="0"/>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$r.InsertXML($xml)
